# The author removed the post in row 657 ("「花を踏みつけることは出来ても
# 春を遅らせることは出来ない」...") entirely from the sheet. Deleting the
# whole row shifts every subsequent row (658-846) up by one, which matches
# the diff (dimension shrinks from C846 to C845, and all rows from 658
# onward are renumbered down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(657).Delete()
